$d = $word.ActiveDocument
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "" -and $i -gt 1 -and $d.Paragraphs.Item($i - 1).Range.Text.Trim() -eq "Q2.2") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Target paragraph (empty paragraph after Q2.2) not found"
}

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:pPr><w:bidi/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:rtl/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>נוכיח שלא ניתן לבצע סכמה במסגרת תיוג רצפים שתאפשר טיפול בישויות מקוננות:</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">תחילה נראה למה </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="cs"/><w:lang w:val="en-US"/></w:rPr><w:t>BILOU</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> לא עובד. עבור המשפט </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>׳׳הועד הפועל של רשות השידור׳׳</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> ישנם 2 ישויות, כל ה</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>׳׳ רשות השידור׳׳</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. שתי הישויות נגמרות במילה "השידור" ולכן צריך תג מיוחד שיבשר ש-2 הישויות נגמרות במילה הזאת והתג </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:lang w:val="en-US"/></w:rPr><w:t>L</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> לא מספיק.</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>אפשר להבין אינטואיטיבית שכדי לפתור את הבעיה צריך שימוש בסוגריים, למשל:</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>׳׳</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">הועד הפועל של </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>רשות השידור</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>׳׳</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">בשביל 2 הסוגריים בסוף, צריך תג מיוחד (למשל </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:lang w:val="en-US"/></w:rPr><w:t>LL</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>, שיהיה כמו 2 סוגרים). אותו דבר אם 2 ישויות שונות היו גם מתחילות באותה המילה.</w:t></w:r></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:rtl/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:bidi/><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>בחזרה לשאלה, נניח בשלילה שאכן יש סכמה כזאת. כמות התגים בסכמה הזאת צריכה להיות סופית, נסמנה ב-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial"/><w:i/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><m:t>M</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> ניצור טקסט</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><w:t>\</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">משפט באורך גדול מזה שמכיל </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><m:t>N</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> ישויות מקוננות (כך ש:</w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><m:t>N&gt;M</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cs="Arial" w:hint="cs"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>) שנגמרות באותה המילה. מאחר ו-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial"/><w:i/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><m:t>N&gt;M</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial" w:hint="cs"/><w:i/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, אין תג שיבשר את סיום </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial"/><w:lang w:val="en-US"/></w:rPr><m:t>N</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math" w:cs="Arial" w:hint="cs"/><w:i/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> הישויות המקוננות, לכן הסכמה לא פותרת את הבעיה.</w:t></w:r></w:p>
'@

$target.Range.InsertXML($xml)
Write-Output "Inserted proof content after Q2.2"
